$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header F1 with same style as E1 (copy format first)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill in time_taken values for rows 2-64
$ws.Range("F2").Value = "2021-10-05 13:39:43.055102"
$ws.Range("F3").Value = "2021-10-05 13:39:43.055115"
$ws.Range("F4").Value = "2021-10-05 13:39:43.055119"
$ws.Range("F5").Value = "2021-10-05 13:39:43.055122"
$ws.Range("F6").Value = "2021-10-05 13:39:43.055125"
$ws.Range("F7").Value = "2021-10-05 13:39:43.055128"
$ws.Range("F8").Value = "2021-10-05 13:39:43.055131"
$ws.Range("F9").Value = "2021-10-05 13:39:43.055134"
$ws.Range("F10").Value = "2021-10-05 13:39:43.055138"
$ws.Range("F11").Value = "2021-10-05 13:39:43.055141"
$ws.Range("F12").Value = "2021-10-05 13:39:43.055144"
$ws.Range("F13").Value = "2021-10-05 13:39:43.055147"
$ws.Range("F14").Value = "2021-10-05 13:39:43.055150"
$ws.Range("F15").Value = "2021-10-05 13:39:43.055153"
$ws.Range("F16").Value = "2021-10-05 13:39:43.055156"
$ws.Range("F17").Value = "2021-10-05 13:39:43.055159"
$ws.Range("F18").Value = "2021-10-05 13:39:43.055162"
$ws.Range("F19").Value = "2021-10-05 13:39:43.055165"
$ws.Range("F20").Value = "2021-10-05 13:39:43.055168"
$ws.Range("F21").Value = "2021-10-05 13:39:43.055171"
$ws.Range("F22").Value = "2021-10-05 13:39:43.055174"
$ws.Range("F23").Value = "2021-10-05 13:39:43.055177"
$ws.Range("F24").Value = "2021-10-05 13:39:43.055180"
$ws.Range("F25").Value = "2021-10-05 13:39:43.055183"
$ws.Range("F26").Value = "2021-10-05 13:39:43.055186"
$ws.Range("F27").Value = "2021-10-05 13:39:43.055189"
$ws.Range("F28").Value = "2021-10-05 13:39:43.055193"
$ws.Range("F29").Value = "2021-10-05 13:39:43.055195"
$ws.Range("F30").Value = "2021-10-05 13:39:43.055198"
$ws.Range("F31").Value = "2021-10-05 13:39:43.055201"
$ws.Range("F32").Value = "2021-10-05 13:39:43.055204"
$ws.Range("F33").Value = "2021-10-05 13:39:43.055207"
$ws.Range("F34").Value = "2021-10-05 13:39:43.055211"
$ws.Range("F35").Value = "2021-10-05 13:39:43.055214"
$ws.Range("F36").Value = "2021-10-05 13:39:43.055217"
$ws.Range("F37").Value = "2021-10-05 13:39:43.055220"
$ws.Range("F38").Value = "2021-10-05 13:39:43.055223"
$ws.Range("F39").Value = "2021-10-05 13:39:43.055256"
$ws.Range("F40").Value = "2021-10-05 13:39:43.055286"
$ws.Range("F41").Value = "2021-10-05 13:39:43.055293"
$ws.Range("F42").Value = "2021-10-05 13:39:43.055297"
$ws.Range("F43").Value = "2021-10-05 13:39:43.055301"
$ws.Range("F44").Value = "2021-10-05 13:39:43.055304"
$ws.Range("F45").Value = "2021-10-05 13:39:43.055307"
$ws.Range("F46").Value = "2021-10-05 13:39:43.055310"
$ws.Range("F47").Value = "2021-10-05 13:39:43.055313"
$ws.Range("F48").Value = "2021-10-05 13:39:43.055316"
$ws.Range("F49").Value = "2021-10-05 13:39:43.055319"
$ws.Range("F50").Value = "2021-10-05 13:39:43.055322"
$ws.Range("F51").Value = "2021-10-05 13:39:43.055325"
$ws.Range("F52").Value = "2021-10-05 13:39:43.055328"
$ws.Range("F53").Value = "2021-10-05 13:39:43.055331"
$ws.Range("F54").Value = "2021-10-05 13:39:43.055335"
$ws.Range("F55").Value = "2021-10-05 13:39:43.055338"
$ws.Range("F56").Value = "2021-10-05 13:39:43.055341"
$ws.Range("F57").Value = "2021-10-05 13:39:43.055344"
$ws.Range("F58").Value = "2021-10-05 13:39:43.055347"
$ws.Range("F59").Value = "2021-10-05 13:39:43.055351"
$ws.Range("F60").Value = "2021-10-05 13:39:43.055353"
$ws.Range("F61").Value = "2021-10-05 13:39:43.055357"
$ws.Range("F62").Value = "2021-10-05 13:39:43.055360"
$ws.Range("F63").Value = "2021-10-05 13:39:43.055363"
$ws.Range("F64").Value = "2021-10-05 13:39:43.055366"
